$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply all cell updates from the crypto price refresh.
# D and E columns hold text values (prices/percentages formatted as
# strings), so force text number format before assigning to avoid
# Excel auto-converting them to numeric values and losing formatting
# (e.g. trailing zeros, thousand-dot separators).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.180.71"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.11%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.714.15"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.30%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "611.35"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.08"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.88%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.590"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.80%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.126"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +6.81%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.08"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.404"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.39%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.46%  "
$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "30.33"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.61%  "
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000209"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +15.80%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.196.72"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.19%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.024.82"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.705.86"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.84"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.98%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "362.37"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.55"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.75%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.25"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.77"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000107"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +12.12%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -4.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.72"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.172"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +4.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.33"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.22"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.45%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "540.32"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.67%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.82"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.67"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.50"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -5.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.437"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.88"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +3.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "163.49"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.01"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.998"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "171.32"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.66%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.94"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.20"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.68%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0619"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.85%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.33"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.75"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.664"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0267"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +4.75%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.46"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0989"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.32%  "
